# Applies scheduled market-data refresh updates to the Jenova Profits workbook.
# For each affected leve row, recompute currentAveragePrice / currentAveragePriceNQ /
# currentAveragePriceHQ / LevePriceNQ / LevePriceHQ / LeveProfitNQ / LeveProfitHQ values.
$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 622.5
$ws.Range("J18").Value = 745
$ws.Range("L18").Value = 745
$ws.Range("N18").Value = -1313
$ws.Range("H40").Value = 7888.3335
$ws.Range("J40").Value = 8110.1665
$ws.Range("L40").Value = 8110.1665
$ws.Range("N40").Value = -8460.166499999999
$ws.Range("H57").Value = 25990
$ws.Range("J57").Value = 8485
$ws.Range("L57").Value = 25455
$ws.Range("N57").Value = -26453
$ws.Range("H98").Value = 2037.5714
$ws.Range("I98").Value = 1463.3055
$ws.Range("J98").Value = 5483.1665
$ws.Range("K98").Value = 1463.3055
$ws.Range("L98").Value = 5483.1665
$ws.Range("M98").Value = 34.69450000000006
$ws.Range("N98").Value = -8479.166499999999
$ws.Range("H116").Value = 8514.591
$ws.Range("I116").Value = 3713.625
$ws.Range("J116").Value = 21317.166
$ws.Range("K116").Value = 3713.625
$ws.Range("L116").Value = 21317.166
$ws.Range("M116").Value = -271.625
$ws.Range("N116").Value = -28201.166
$ws.Range("H122").Value = 2037.5714
$ws.Range("I122").Value = 1463.3055
$ws.Range("J122").Value = 5483.1665
$ws.Range("K122").Value = 4389.916499999999
$ws.Range("L122").Value = 16449.4995
$ws.Range("M122").Value = -1939.916499999999
$ws.Range("N122").Value = -21349.4995
$ws.Range("H125").Value = 7940253
$ws.Range("J125").Value = 10105149
$ws.Range("L125").Value = 90946341
$ws.Range("N125").Value = -90951261
$ws.Range("H132").Value = 2759.2163
$ws.Range("I132").Value = 2633.3
$ws.Range("K132").Value = 7899.900000000001
$ws.Range("M132").Value = -5369.900000000001
$ws.Range("H137").Value = 772429.25
$ws.Range("I137").Value = 627578.9
$ws.Range("K137").Value = 1882736.7
$ws.Range("M137").Value = -1880186.7
$ws.Range("H138").Value = 4326.269
$ws.Range("J138").Value = 5598.4307
$ws.Range("L138").Value = 16795.2921
$ws.Range("N138").Value = -27075.2921

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H59").Value = 86714.25
$ws.Range("J59").Value = 86714.25
$ws.Range("L59").Value = 86714.25
$ws.Range("N59").Value = -88408.25
$ws.Range("H86").Value = 896216.3
$ws.Range("I86").Value = 1215900
$ws.Range("K86").Value = 1215900
$ws.Range("M86").Value = -1214777
$ws.Range("H89").Value = 896216.3
$ws.Range("I89").Value = 1215900
$ws.Range("K89").Value = 6079500
$ws.Range("M89").Value = -6073884
$ws.Range("H99").Value = 1207.4445
$ws.Range("J99").Value = 1025
$ws.Range("L99").Value = 1025
$ws.Range("N99").Value = -4021
$ws.Range("H139").Value = 77810
$ws.Range("J139").Value = 77810
$ws.Range("L139").Value = 77810
$ws.Range("N139").Value = -88090

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 112.05556
$ws.Range("I7").Value = 40.46154
$ws.Range("J7").Value = 298.2
$ws.Range("K7").Value = 40.46154
$ws.Range("L7").Value = 298.2
$ws.Range("M7").Value = 72.53846
$ws.Range("N7").Value = -524.2
$ws.Range("H31").Value = 80684.84
$ws.Range("I31").Value = 0
$ws.Range("J31").Value = 80684.84
$ws.Range("K31").Value = 0
$ws.Range("L31").Value = 80684.84
$ws.Range("M31").ClearContents()
$ws.Range("N31").Value = -81274.84
$ws.Range("H34").Value = 80684.84
$ws.Range("I34").Value = 0
$ws.Range("J34").Value = 80684.84
$ws.Range("K34").Value = 0
$ws.Range("L34").Value = 80684.84
$ws.Range("M34").ClearContents()
$ws.Range("N34").Value = -81088.84
$ws.Range("H41").Value = 31690
$ws.Range("I41").Value = 28019.666
$ws.Range("J41").Value = 32791.1
$ws.Range("K41").Value = 28019.666
$ws.Range("L41").Value = 32791.1
$ws.Range("M41").Value = -27591.666
$ws.Range("N41").Value = -33647.1
$ws.Range("H62").Value = 3960.5557
$ws.Range("J62").Value = 5679
$ws.Range("L62").Value = 5679
$ws.Range("N62").Value = -6927
$ws.Range("H65").Value = 3960.5557
$ws.Range("J65").Value = 5679
$ws.Range("L65").Value = 28395
$ws.Range("N65").Value = -34635
$ws.Range("H122").Value = 4561.7856
$ws.Range("I122").Value = 2564.625
$ws.Range("J122").Value = 7224.6665
$ws.Range("K122").Value = 7693.875
$ws.Range("L122").Value = 21673.9995
$ws.Range("M122").Value = -5243.875
$ws.Range("N122").Value = -26573.9995
$ws.Range("H139").Value = 74875
$ws.Range("J139").Value = 74875
$ws.Range("L139").Value = 74875
$ws.Range("N139").Value = -85155

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H75").Value = 2830
$ws.Range("I75").Value = 1495
$ws.Range("J75").Value = 3497.5
$ws.Range("K75").Value = 4485
$ws.Range("L75").Value = 10492.5
$ws.Range("M75").Value = -3487
$ws.Range("N75").Value = -12488.5
$ws.Range("H78").Value = 2830
$ws.Range("I78").Value = 1495
$ws.Range("J78").Value = 3497.5
$ws.Range("K78").Value = 13455
$ws.Range("L78").Value = 31477.5
$ws.Range("M78").Value = -8463
$ws.Range("N78").Value = -41461.5
$ws.Range("H134").Value = 979.9091
$ws.Range("I134").Value = 979.9091
$ws.Range("K134").Value = 2939.7273
$ws.Range("M134").Value = 2130.2727
$ws.Range("H140").Value = 2143.4412
$ws.Range("I140").Value = 1429.2333
$ws.Range("K140").Value = 4287.699900000001
$ws.Range("M140").Value = 892.3000999999995
$ws.Range("H141").Value = 6460.273
$ws.Range("I141").Value = 6460.273
$ws.Range("K141").Value = 19380.819
$ws.Range("M141").Value = -14200.819

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 5060.905
$ws.Range("I61").Value = 4774
$ws.Range("K61").Value = 4774
$ws.Range("M61").Value = -4572
$ws.Range("H68").Value = 2893
$ws.Range("J68").Value = 3249.4
$ws.Range("L68").Value = 3249.4
$ws.Range("N68").Value = -4747.4
$ws.Range("H71").Value = 2893
$ws.Range("J71").Value = 3249.4
$ws.Range("L71").Value = 16247
$ws.Range("N71").Value = -23735
$ws.Range("H113").Value = 5060.905
$ws.Range("I113").Value = 4774
$ws.Range("K113").Value = 4774
$ws.Range("M113").Value = -2604
$ws.Range("H122").Value = 502248.44
$ws.Range("I122").Value = 2220.2222
$ws.Range("K122").Value = 6660.6666
$ws.Range("M122").Value = -4210.6666
$ws.Range("H132").Value = 5431.9287
$ws.Range("I132").Value = 2918
$ws.Range("J132").Value = 6117.5454
$ws.Range("K132").Value = 8754
$ws.Range("L132").Value = 18352.6362
$ws.Range("M132").Value = -6224
$ws.Range("N132").Value = -23412.6362
$ws.Range("H136").Value = 1118410.6
$ws.Range("I136").Value = 1544806.4
$ws.Range("K136").Value = 4634419.199999999
$ws.Range("M136").Value = -4631869.199999999
$ws.Range("H139").Value = 53999.832
$ws.Range("J139").Value = 55000
$ws.Range("L139").Value = 55000
$ws.Range("N139").Value = -65280

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 28573978
$ws.Range("I122").Value = 31252086
$ws.Range("K122").Value = 93756258
$ws.Range("M122").Value = -93753808
$ws.Range("H126").Value = 968
$ws.Range("I126").Value = 1046
$ws.Range("K126").Value = 3138
$ws.Range("M126").Value = -668
$ws.Range("H132").Value = 17741.953
$ws.Range("I132").Value = 1555
$ws.Range("J132").Value = 105151.5
$ws.Range("K132").Value = 4665
$ws.Range("L132").Value = 315454.5
$ws.Range("M132").Value = -2135
$ws.Range("N132").Value = -320514.5
$ws.Range("H136").Value = 6908634.5
$ws.Range("J136").Value = 146509.92
$ws.Range("L136").Value = 439529.76
$ws.Range("N136").Value = -444629.76
$ws.Range("H138").Value = 82963.5
$ws.Range("J138").Value = 82963.5
$ws.Range("L138").Value = 82963.5
$ws.Range("N138").Value = -93243.5
